# Add IR GeneralInformation test script data: new columns AB:AG with
# ArrestDetails / WorkplaceViolence / TypeofAccident / RequiredField /
# MandatoryonSave / MandatoryonClose, plus the corresponding row-2 test
# values, and fill in a few previously-blank row-2 cells (Title, a new
# "No" dropdown cell, ReportingDate "NULL", plus the boolean notification
# flags that were missing).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New header row (row 1), columns AB:AG -------------------------------
$ws.Range("AB1").Value = "ArrestDetails"
$ws.Range("AC1").Value = "WorkplaceViolence"
$ws.Range("AD1").Value = "TypeofAccident "
$ws.Range("AE1").Value = "RequiredField"
$ws.Range("AF1").Value = "MandatoryonSave"
$ws.Range("AG1").Value = "MandatoryonClose "

# --- Row 2: fill previously-empty cells (order matches the original
#     authoring order so new shared-string indices line up) ----------------
$ws.Range("AB2").Value = "No"
$ws.Range("AC2").Value = "No"
$ws.Range("AD2").Value = "Posture"
$ws.Range("AE2").Value = "Hitech"
$ws.Range("AF2").Value = "Save"
$ws.Range("AG2").Value = "Close"

$ws.Range("G2").Value = "Automation test new Incident Report "
$ws.Range("M2").Value = "NULL"

$ws.Range("N2").Value = 1
$ws.Range("O2").Value = 28

$ws.Range("U2").Value = $true
$ws.Range("V2").Value = $true
$ws.Range("W2").Value = $true
$ws.Range("Y2").Value = $true
$ws.Range("Z2").Value = $true
$ws.Range("AA2").Value = $true

# --- New column widths (AB:AG) --------------------------------------------
$ws.Columns.Item(28).ColumnWidth = 14.5   # AB
$ws.Columns.Item(29).ColumnWidth = 15.5   # AC
$ws.Columns.Item(30).ColumnWidth = 16     # AD
$ws.Columns.Item(32).ColumnWidth = 17.5   # AF
$ws.Columns.Item(33).ColumnWidth = 18.5   # AG

# --- Alignment: left-align row 2 and the previously centre-only cells -----
$ws.Range("A2").HorizontalAlignment = -4131
$ws.Range("B2:AG2").HorizontalAlignment = -4131
$ws.Range("AB1:AG1").EntireColumn.VerticalAlignment = -4108

# --- Selection as left by the author on save -------------------------------
[void]$ws.Range("A2").Select()
